# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
#
# This script re-sorts a handful of match rows in the
# "Germany Verbandsliga" sheet:
#   - rows 38/39/40 are cyclically rotated (new38<-old40, new39<-old38, new40<-old39)
#   - rows 132/133 are swapped (new132<-old133, new133<-old132)
# Columns A (row #), C, D and E are identical across every affected row
# (same Div / Div Original Name / Date), so only columns B and F..AC
# actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Verbandsliga")

# --- rows 38 / 39 / 40 : 3-way cyclic rotation -----------------------
$v38 = $ws.Range("B38:AC38").Value()
$v39 = $ws.Range("B39:AC39").Value()
$v40 = $ws.Range("B40:AC40").Value()

$ws.Range("B38:AC38").Value = $v40
$ws.Range("B39:AC39").Value = $v38
$ws.Range("B40:AC40").Value = $v39

# --- rows 132 / 133 : simple swap -------------------------------------
$v132 = $ws.Range("B132:AC132").Value()
$v133 = $ws.Range("B133:AC133").Value()

$ws.Range("B132:AC132").Value = $v133
$ws.Range("B133:AC133").Value = $v132

Write-Host "Done"
